$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing columns (X,Y) that no longer exist in the new layout
$ws.Range("X1:Y4").Clear()

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Servoname"
$ws.Range("B1").Value = "length"
$ws.Range("C1").Value = "height"
$ws.Range("D1").Value = "width"
$ws.Range("E1").Value = "shaft_R"
$ws.Range("F1").Value = "shaft_offs"
$ws.Range("G1").Value = "connect_R"
$ws.Range("H1").Value = "connect_top_R"
$ws.Range("I1").Value = "connect_top_H"
$ws.Range("J1").Value = "connect_screw_R"
$ws.Range("K1").Value = "connect_screw_circle_R"
$ws.Range("L1").Value = "connect_screw_depth"
$ws.Range("M1").Value = "connect_screw_Num"
$ws.Range("N1").Value = "connect_dis"
$ws.Range("O1").Value = "PL_cable_gap_hor"
$ws.Range("P1").Value = "PL_cable_gap_ver"
$ws.Range("Q1").Value = "mounting_screw_R"
$ws.Range("R1").Value = "mounting_screw_depth"
$ws.Range("S1").Value = "screw_mount_y"
$ws.Range("T1").Value = "screw_mount_x"
$ws.Range("U1").Value = "screw_mount_z"
$ws.Range("V1").Value = "cable_gap_width"
$ws.Range("W1").Value = "cable_gap"

# ---- Row 2 (sm40bl) ----
$ws.Range("A2").Value = "sm40bl"
$ws.Range("B2").Value = 46.5
$ws.Range("C2").Value = 34
$ws.Range("D2").Value = 28.5
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 11.25
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 2.5
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = "[-22.25,-2.5;22.25,-2.5;22.25,12.5;-22.25,12.5]"
$ws.Range("P2").Value = "[-22.25,-15.5;22.25,-15.5;22.25,-3;-22.25,-3]"
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = "[-8,6;8,6;-8,-6;8,-6]"
$ws.Range("T2").Value = "[6,-16;-6,-16;6,8]"
$ws.Range("U2").Value = "NaN"
$ws.Range("V2").Value = 20
$ws.Range("W2").Value = 8

# ---- Row 3 (sm85bl) ----
$ws.Range("A3").Value = "sm85bl"
$ws.Range("B3").Value = 62
$ws.Range("C3").Value = 47
$ws.Range("D3").Value = 34
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 10.5
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 8
$ws.Range("N3").Value = 55.6
$ws.Range("O3").Value = "[-25,-7.5;25,-7.5;25,7.5;-25,7.5]"
$ws.Range("P3").Value = "[-27,-15;27,-15;27,5;-27,5]"
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = "NaN"
$ws.Range("T3").Value = "NaN"
$ws.Range("U3").Value = "[14,0;-14,0;14,28;-14,28;14,-28;-14,-28]"
$ws.Range("V3").Value = 20
$ws.Range("W3").Value = 8

# ---- Row 4 (sm120bl) ----
$ws.Range("A4").Value = "sm120bl"
$ws.Range("B4").Value = 78
$ws.Range("C4").Value = 61.5
$ws.Range("D4").Value = 43
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 21.5
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 18
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 12.5
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 72.6
$ws.Range("O4").Value = "[-15,40;15,40;15,60;-15,60]"
$ws.Range("P4").Value = "NaN"
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = "NaN"
$ws.Range("T4").Value = "NaN"
$ws.Range("U4").Value = "[17.5,0;-17.5,0;17.5,35;-17.5,35;17.5,-35;-17.5,-35]"
$ws.Range("V4").Value = 20
$ws.Range("W4").Value = 8

# ---- Style bookkeeping: the one "applyNumberFormat" style moves from V2 to L2 ----
$ws.Range("L2").NumberFormat = "General"
$ws.Range("V2").ClearFormats()
$ws.Range("L2").Value = 2.5
$ws.Range("V2").Value = 20

# ---- Column widths (approximate best-fit sizing for the new layout) ----
$ws.Columns("F").ColumnWidth = 9.85546875
$ws.Columns("G").ColumnWidth = 10.140625
$ws.Columns("H").ColumnWidth = 14.140625
$ws.Columns("I").ColumnWidth = 14.28515625
$ws.Columns("J").ColumnWidth = 16.42578125
$ws.Columns("K").ColumnWidth = 22.28515625
$ws.Columns("L").ColumnWidth = 20.5703125
$ws.Columns("M").ColumnWidth = 19.5703125
$ws.Columns("N").ColumnWidth = 11.5703125
$ws.Columns("O").ColumnWidth = 36.140625
$ws.Columns("P").ColumnWidth = 38.85546875
$ws.Columns("Q").ColumnWidth = 18
$ws.Columns("R").ColumnWidth = 19.42578125
$ws.Columns("S").ColumnWidth = 17.42578125
$ws.Columns("T").ColumnWidth = 16.28515625
$ws.Columns("U").ColumnWidth = 44.5703125
$ws.Columns("V").ColumnWidth = 16.140625

# ---- Selection / view state ----
$ws.Range("K12").Select()

# ---- Workbook window geometry ----
$excel.Width = 1788
$excel.Height = 972
